# piv-opt.pptx edit: "Change PQO to GPORCA in piv-opt.png graphic (#1886)"
#
# 1. Bump the cached datetimeFigureOut placeholder text (6/6/2014 -> 2/2/2017)
#    on the slide master and every slide layout.
# 2. Rename the "Pivotal Query Optimizer" box to "GPORCA" and shrink/reposition
#    its text box to fit the shorter single-line label.
# 3. Consolidate the two "Legacy " + "Optimizer" runs into a single
#    "Legacy Optimizer" run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder text: Slide Master + all Custom Layouts
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "6/6/2014") {
                $shp.TextFrame.TextRange.Text = "2/2/2017"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i)
}

# ---------------------------------------------------------------------------
# 2 & 3. Slide 1 shapes
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# Locate the two group shapes that hold the optimizer labels by walking all
# top-level shapes and inspecting their GroupItems (robust to index drift).
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $top = $s.Shapes.Item($i)
    if ($top.Type -eq 6) {
        for ($k = 1; $k -le $top.GroupItems.Count; $k++) {
            $inner = $top.GroupItems.Item($k)
            if (-not $inner.HasTextFrame) { continue }

            $txt = $inner.TextFrame.TextRange.Text

            if ($txt -eq "Pivotal Query Optimizer") {
                # Rename + shrink/reposition the text box (autofit shrinks the
                # height to match the new single-line text; only Top needs to
                # be moved explicitly to its new EMU-exact position).
                $inner.TextFrame.TextRange.Text = "GPORCA"
                $inner.Top = 199.99216535433072
            }
            elseif ($txt -eq "Legacy Optimizer") {
                # Force a full-run rebuild so the two runs ("Legacy " and
                # "Optimizer") collapse into a single run, matching the diff.
                $inner.TextFrame.TextRange.Text = "Legacy OptimizerX"
                $inner.TextFrame.TextRange.Text = "Legacy Optimizer"
            }
        }
    }
}
